# Applies the "LTSD Parameters" table addition to the APA analysis sheet.
#
# Final layout (rows/columns) after the edit:
#   K2 = "LTSD Parameters"                          (was "Test")
#   K3 = "Right"                      M3 = "Left"
#   K4 = "Threshols"  L4 = "Win"      M4 = "Threshold"  N4 = "Win"
#   K5 = "5.5"        L5 = "300.0"    M5 = "5.0"         N5 = "300.0"
#
# The existing filename list (A5:A14), the GaitWatch file reference (D5) and
# the duration value (E5) stay exactly where they were - only new columns
# (K:N) and new rows (3 and 4) are introduced; no existing rows move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that *looks* numeric ("5.5", "300.0", ...) as genuine
# text (matching the source workbook, where these are shared strings, not
# numbers) by round-tripping it through a text-formatted scratch cell and a
# values-only paste, then restoring the scratch cell to its original state.
function Set-TextValue {
    param($range, [string]$text)

    $scratch = $ws.Range("ZZ9000")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null   # xlPasteValues
    $scratch.Clear() | Out-Null
}

# --- Row 2: replace "Test" in K2 with the new header ---
$ws.Range("K2").Value = "LTSD Parameters"

# --- Row 3: Right / Left headers ---
$ws.Range("K3").Value = "Right"
$ws.Range("M3").Value = "Left"

# --- Row 4: sub headers ---
$ws.Range("K4").Value = "Threshols"
$ws.Range("L4").Value = "Win"
$ws.Range("M4").Value = "Threshold"
$ws.Range("N4").Value = "Win"

# --- Row 5: parameter values (kept as text, e.g. "5.0" not 5) ---
# Written in this order so new shared-string entries land in the same
# sequence as the source workbook (5.0, 5.5, 300.0).
Set-TextValue $ws.Range("M5") "5.0"
Set-TextValue $ws.Range("K5") "5.5"
Set-TextValue $ws.Range("L5") "300.0"
Set-TextValue $ws.Range("N5") "300.0"

# Update the selection to match the post-edit state (active cell N5).
$ws.Range("N5").Select() | Out-Null

$wb.Save()
